# Update the single-column results table in place.
# Each entry below is a (row, newValue) pair: the table has one column,
# so we address cells as Cell(row, 1) and overwrite the cell text wholesale.
# This also collapses the multi-run/tab-separated "summary" rows (44-46)
# down to a single simple value, matching the target content.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$updates = @(
    @{Row = 1;  Value = "0M"},
    @{Row = 2;  Value = "0M"},
    @{Row = 3;  Value = "0M"},
    @{Row = 4;  Value = "111"},
    @{Row = 5;  Value = "0.00004"},
    @{Row = 7;  Value = "0.00018"},
    @{Row = 9;  Value = "0.00038"},
    @{Row = 10; Value = "0.00042"},
    @{Row = 11; Value = "0.00044"},
    @{Row = 12; Value = "0.02358"},
    @{Row = 44; Value = "99.96"},
    @{Row = 45; Value = "0.02"},
    @{Row = 46; Value = "55"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, 1)
    $cell.Range.Text = $u.Value
}
